$d = $word.ActiveDocument

$replacements = @(
    @("2024-03-09 Saturday", "2024-03-10 Sunday"),
    @("241×5=", "264×2="),
    @("533×8=", "978×8="),
    @("827×3=", "398×8="),
    @("723×4=", "756×6="),
    @("547×3=", "385×8="),
    @("523×7=", "897×6="),
    @("498×7=", "968×4="),
    @("157×5=", "198×9="),
    @("660×5=", "679×4="),
    @("192×4=", "566×2="),
    @("178×8=", "529×9="),
    @("479×9=", "434×8="),
    @("446×8=", "423×4="),
    @("219×4=", "183×9="),
    @("743×4=", "448×5="),
    @("759×6=", "102×9="),
    @("472×2=", "143×6="),
    @("588×8=", "272×7="),
    @("517×7=", "664×5="),
    @("624×9=", "299×5="),
    @("712×9=", "926×6="),
    @("746×8=", "138×4="),
    @("270×9=", "163×3="),
    @("612×9=", "423×7="),
    @("511×8=", "813×7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
